# Update "想去人数" (F column) figures to the latest scraped counts
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 778
$ws.Range("F7").Value = 406
$ws.Range("F8").Value = 853
$ws.Range("F10").Value = 6996
$ws.Range("F11").Value = 1928
$ws.Range("F12").Value = 5097
$ws.Range("F13").Value = 499
$ws.Range("F15").Value = 6632
$ws.Range("F16").Value = 8225
$ws.Range("F20").Value = 4230
$ws.Range("F21").Value = 630
$ws.Range("F22").Value = 114
$ws.Range("F23").Value = 0
$ws.Range("F26").Value = 1138
$ws.Range("F27").Value = 53
$ws.Range("F29").Value = 644
$ws.Range("F30").Value = 789
$ws.Range("F31").Value = 1790
$ws.Range("F33").Value = 2106
$ws.Range("F34").Value = 269
$ws.Range("F35").Value = 82
$ws.Range("F36").Value = 1358
$ws.Range("F38").Value = 754
$ws.Range("F39").Value = 363
$ws.Range("F40").Value = 2835
$ws.Range("F41").Value = 3885
$ws.Range("F45").Value = 483
$ws.Range("F46").Value = 8
$ws.Range("F47").Value = 844
$ws.Range("F48").Value = 135
$ws.Range("F49").Value = 4020

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 1339
$ws.Range("F31").Value = 96

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 4859

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 4859
$ws.Range("F5").Value = 1339
$ws.Range("F11").Value = 778
$ws.Range("F12").Value = 406
$ws.Range("F13").Value = 853
$ws.Range("F17").Value = 5097
$ws.Range("F18").Value = 6632
$ws.Range("F19").Value = 6632
$ws.Range("F24").Value = 4230
$ws.Range("F25").Value = 630
$ws.Range("F26").Value = 114
$ws.Range("F29").Value = 1138
$ws.Range("F31").Value = 644
$ws.Range("F32").Value = 789
$ws.Range("F33").Value = 1790
$ws.Range("F35").Value = 2106
$ws.Range("F40").Value = 754
$ws.Range("F42").Value = 363
$ws.Range("F43").Value = 96
$ws.Range("F44").Value = 3885
$ws.Range("F47").Value = 844
$ws.Range("F48").Value = 135
$ws.Range("F50").Value = 4020

